$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4,1,2,2,1,1,1,1,1,4,2,2,1),
    @(5,1,2,2,1,1,1,1,1,5,2,2,1),
    @(6,1,2,2,1,1,1,1,1,6,2,2,1),
    @(7,1,3,2,1,1,1,1,1,7,2,2,1),
    @(8,1,3,2,1,1,1,1,1,8,2,2,1),
    @(9,1,3,2,1,1,1,1,1,9,2,2,1)
)

$startRow = 5
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}

$ws.Range("C11").Select()
